# AssetAddress.xlsx — UI MVVM content update
# Replace the "Cube" row's path with the new Resource-qualified path and add
# two new asset rows ("sword" and "broom"), then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same, but re-assert it for safety.
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "path"

# Row 2: existing "Cube" entry now points at Assets/Resource/Cube.prefab.
$ws.Range("A2").Value = "Cube"
$ws.Range("B2").Value = "Assets/Resource/Cube.prefab"

# Row 3: new "sword" entry (previously held plain numeric placeholder data).
$ws.Range("A3").Value = "sword"
$ws.Range("B3").Value = "Assets/Resource/sword.png"

# Row 4: new "broom" entry (previously held plain numeric placeholder data).
$ws.Range("A4").Value = "broom"
$ws.Range("B4").Value = "Assets/Resource/broom.png"

# Move the sheet's active selection from E6 to G6.
$ws.Range("G6").Select()
